$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Cat"
$ws.Range("B6").Value = "incomplete"
$ws.Range("C6").Value = "2025-01-05 21:33:59.848542"
$ws.Range("D6").Value = "N/A"

$ws.Range("A7").Value = "at"
$ws.Range("B7").Value = "incomplete"
$ws.Range("C7").Value = "2025-01-05 21:38:36.142019"
$ws.Range("D7").Value = "N/A"
